$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows 129-131 ---
$ws.Range("Z129").Value = 13.3529615
$ws.Range("AR129").Value = 7.5893467
$ws.Range("Z130").Value = 13.1781499
$ws.Range("AR130").Value = 10.2485045
$ws.Range("BE130").Value = 14.9564713
$ws.Range("F131").Value = 16.5461133
$ws.Range("G131").Value = 11.2959453
$ws.Range("L131").Value = 11.7876606
$ws.Range("O131").Value = 8.5590189
$ws.Range("S131").Value = 18.8809763
$ws.Range("Z131").Value = 13.1246331
$ws.Range("AT131").Value = 15.0204726
$ws.Range("AV131").Value = 13.9230571
$ws.Range("BB131").Value = 12.4514466
$ws.Range("BD131").Value = 13.9446283
$ws.Range("BE131").Value = 15.1876081

# --- Fill in full data for existing label rows 132-133, and new rows 134-138 ---
# Row 132
$ws.Range("B132").Value = 14.2078255
$ws.Range("C132").Value = 19.301172
$ws.Range("D132").Value = 17.8812427
$ws.Range("F132").Value = 17.1699225
$ws.Range("G132").Value = 11.3159436
$ws.Range("H132").Value = 13.7565069
$ws.Range("I132").Value = 13.6173913
$ws.Range("J132").Value = 14.8976982
$ws.Range("K132").Value = 13.5141294
$ws.Range("L132").Value = 11.529806
$ws.Range("M132").Value = 15.5494627
$ws.Range("O132").Value = 7.5867426
$ws.Range("P132").Value = 20.1670333
$ws.Range("Q132").Value = 14.7044281
$ws.Range("R132").Value = 14.3339816
$ws.Range("S132").Value = 18.7015281
$ws.Range("T132").Value = 14.3837301
$ws.Range("U132").Value = 15.345495
$ws.Range("V132").Value = 18.4882367
$ws.Range("W132").Value = 14.6836525
$ws.Range("X132").Value = 14.9465457
$ws.Range("Y132").Value = 10.0928481
$ws.Range("Z132").Value = 12.7157998
$ws.Range("AA132").Value = 15.3844612
$ws.Range("AB132").Value = 14.1121981
$ws.Range("AD132").Value = 20.6459457
$ws.Range("AE132").Value = 10.5509873
$ws.Range("AF132").Value = 14.7784578
$ws.Range("AG132").Value = 17.6029318
$ws.Range("AH132").Value = 20.7043035
$ws.Range("AI132").Value = 12.0955203
$ws.Range("AJ132").Value = 14.3063212
$ws.Range("AK132").Value = 13.9300282
$ws.Range("AL132").Value = 13.4878516
$ws.Range("AM132").Value = 13.3802825
$ws.Range("AN132").Value = 13.567387
$ws.Range("AO132").Value = 14.2578272
$ws.Range("AP132").Value = 12.0109486
$ws.Range("AQ132").Value = 11.7041099
$ws.Range("AS132").Value = 13.3623027
$ws.Range("AT132").Value = 15.6660148
$ws.Range("AU132").Value = 21.0719259
$ws.Range("AV132").Value = 13.9482109
$ws.Range("AW132").Value = 15.2458642
$ws.Range("AX132").Value = 17.770569
$ws.Range("AY132").Value = 14.3638439
$ws.Range("BA132").Value = 8.7316102
$ws.Range("BB132").Value = 12.4410437
$ws.Range("BC132").Value = 13.3063874
$ws.Range("BD132").Value = 14.3944167
$ws.Range("BE132").Value = 14.387196

# Row 133
$ws.Range("B133").Value = 15.4071661
$ws.Range("C133").Value = 19.97973
$ws.Range("D133").Value = 17.5865215
$ws.Range("F133").Value = 17.4152647
$ws.Range("G133").Value = 11.2590772
$ws.Range("H133").Value = 13.3875337
$ws.Range("I133").Value = 13.2490803
$ws.Range("J133").Value = 15.3795812
$ws.Range("K133").Value = 13.3096927
$ws.Range("L133").Value = 11.589119
$ws.Range("M133").Value = 15.1490648
$ws.Range("O133").Value = 7.8556807
$ws.Range("P133").Value = 20.5131229
$ws.Range("Q133").Value = 13.8497096
$ws.Range("R133").Value = 14.0785931
$ws.Range("S133").Value = 18.79548
$ws.Range("T133").Value = 14.3156073
$ws.Range("U133").Value = 15.0097051
$ws.Range("V133").Value = 18.6318984
$ws.Range("W133").Value = 14.3606388
$ws.Range("X133").Value = 14.2663854
$ws.Range("Y133").Value = 10.0711482
$ws.Range("Z133").Value = 12.4796205
$ws.Range("AA133").Value = 15.3820211
$ws.Range("AB133").Value = 13.6507954
$ws.Range("AD133").Value = 21.050917
$ws.Range("AE133").Value = 10.0035525
$ws.Range("AF133").Value = 14.7180546
$ws.Range("AG133").Value = 17.5037305
$ws.Range("AH133").Value = 20.8677006
$ws.Range("AI133").Value = 11.6119483
$ws.Range("AJ133").Value = 14.2969908
$ws.Range("AK133").Value = 13.8412321
$ws.Range("AL133").Value = 12.6325155
$ws.Range("AM133").Value = 13.0372738
$ws.Range("AN133").Value = 13.647672
$ws.Range("AO133").Value = 14.0221268
$ws.Range("AP133").Value = 11.6657081
$ws.Range("AQ133").Value = 11.7214128
$ws.Range("AS133").Value = 12.9219127
$ws.Range("AT133").Value = 16.1878248
$ws.Range("AU133").Value = 20.4633716
$ws.Range("AV133").Value = 13.9762194
$ws.Range("AW133").Value = 15.2708939
$ws.Range("AX133").Value = 18.1371273
$ws.Range("AY133").Value = 14.0725748
$ws.Range("BA133").Value = 8.659206599999999
$ws.Range("BB133").Value = 12.4906953
$ws.Range("BC133").Value = 13.5647435
$ws.Range("BD133").Value = 13.6116076
$ws.Range("BE133").Value = 14.8345475

# Row 134
$ws.Range("A134").Value = "12 06 2020"
$ws.Range("B134").Value = 15.300727
$ws.Range("C134").Value = 19.5356138
$ws.Range("D134").Value = 17.6595259
$ws.Range("F134").Value = 17.6349544
$ws.Range("G134").Value = 11.3925637
$ws.Range("H134").Value = 12.9923334
$ws.Range("I134").Value = 12.9037963
$ws.Range("J134").Value = 14.083558
$ws.Range("K134").Value = 12.7772643
$ws.Range("L134").Value = 11.7384777
$ws.Range("M134").Value = 15.0601687
$ws.Range("O134").Value = 8.299899699999999
$ws.Range("P134").Value = 20.070983
$ws.Range("Q134").Value = 14.6148805
$ws.Range("R134").Value = 13.7458832
$ws.Range("S134").Value = 18.7859837
$ws.Range("T134").Value = 14.2674471
$ws.Range("U134").Value = 14.655373
$ws.Range("V134").Value = 18.7587594
$ws.Range("W134").Value = 14.4270854
$ws.Range("X134").Value = 13.9680509
$ws.Range("Y134").Value = 9.256577399999999
$ws.Range("Z134").Value = 12.014159
$ws.Range("AA134").Value = 15.0357303
$ws.Range("AB134").Value = 14.2304148
$ws.Range("AD134").Value = 21.0802879
$ws.Range("AE134").Value = 9.5368245
$ws.Range("AF134").Value = 14.629002
$ws.Range("AG134").Value = 17.1493065
$ws.Range("AH134").Value = 21.3465108
$ws.Range("AI134").Value = 11.6519174
$ws.Range("AJ134").Value = 13.9783621
$ws.Range("AK134").Value = 14.3683957
$ws.Range("AL134").Value = 12.9616137
$ws.Range("AM134").Value = 12.9308176
$ws.Range("AN134").Value = 13.6416097
$ws.Range("AO134").Value = 14.8479413
$ws.Range("AP134").Value = 11.7670466
$ws.Range("AQ134").Value = 11.6616173
$ws.Range("AS134").Value = 12.9244698
$ws.Range("AT134").Value = 16.559976
$ws.Range("AU134").Value = 20.8628418
$ws.Range("AV134").Value = 13.9550334
$ws.Range("AW134").Value = 15.7336092
$ws.Range("AX134").Value = 18.4351064
$ws.Range("AY134").Value = 13.8164021
$ws.Range("BA134").Value = 8.528719000000001
$ws.Range("BB134").Value = 12.6428402
$ws.Range("BC134").Value = 13.4217538
$ws.Range("BD134").Value = 14.248871
$ws.Range("BE134").Value = 14.6055422

# Row 135
$ws.Range("A135").Value = "13 06 2020"
$ws.Range("B135").Value = 15.1774398
$ws.Range("C135").Value = 19.4539491
$ws.Range("D135").Value = 18.0101713
$ws.Range("F135").Value = 17.9183935
$ws.Range("G135").Value = 11.4129938
$ws.Range("H135").Value = 12.7335787
$ws.Range("I135").Value = 12.680776
$ws.Range("J135").Value = 14.5371578
$ws.Range("K135").Value = 12.8865979
$ws.Range("L135").Value = 11.7591432
$ws.Range("M135").Value = 15.3592905
$ws.Range("O135").Value = 8.100558700000001
$ws.Range("P135").Value = 20.0351179
$ws.Range("Q135").Value = 14.4814364
$ws.Range("R135").Value = 13.4698227
$ws.Range("S135").Value = 18.6605114
$ws.Range("T135").Value = 13.7331304
$ws.Range("U135").Value = 14.426606
$ws.Range("V135").Value = 18.6317649
$ws.Range("W135").Value = 13.9226582
$ws.Range("X135").Value = 13.9851683
$ws.Range("Y135").Value = 9.076433099999999
$ws.Range("Z135").Value = 11.7842875
$ws.Range("AA135").Value = 14.9135091
$ws.Range("AB135").Value = 14.1342852
$ws.Range("AD135").Value = 21.0962135
$ws.Range("AE135").Value = 9.824619999999999
$ws.Range("AF135").Value = 14.9885001
$ws.Range("AG135").Value = 16.6291366
$ws.Range("AH135").Value = 20.5507839
$ws.Range("AI135").Value = 11.5291262
$ws.Range("AJ135").Value = 13.7300075
$ws.Range("AK135").Value = 14.0712605
$ws.Range("AL135").Value = 13.4167147
$ws.Range("AM135").Value = 12.7142376
$ws.Range("AN135").Value = 13.6848956
$ws.Range("AO135").Value = 14.6589586
$ws.Range("AP135").Value = 11.9241602
$ws.Range("AQ135").Value = 11.4898137
$ws.Range("AS135").Value = 11.9037145
$ws.Range("AT135").Value = 16.960511
$ws.Range("AU135").Value = 20.3462191
$ws.Range("AV135").Value = 14.019408
$ws.Range("AW135").Value = 15.9887428
$ws.Range("AX135").Value = 18.4734274
$ws.Range("AY135").Value = 13.6750689
$ws.Range("BA135").Value = 8.146070699999999
$ws.Range("BB135").Value = 12.8200266
$ws.Range("BC135").Value = 13.3065282
$ws.Range("BD135").Value = 13.661591
$ws.Range("BE135").Value = 13.6279081

# Row 136
$ws.Range("A136").Value = "14 06 2020"
$ws.Range("B136").Value = 14.9645809
$ws.Range("C136").Value = 19.8181272
$ws.Range("D136").Value = 18.4630463
$ws.Range("F136").Value = 18.2413088
$ws.Range("G136").Value = 11.4578952
$ws.Range("H136").Value = 12.5188231
$ws.Range("I136").Value = 12.1695906
$ws.Range("J136").Value = 14.0977444
$ws.Range("K136").Value = 12.2875517
$ws.Range("L136").Value = 11.9147423
$ws.Range("M136").Value = 15.2651201
$ws.Range("O136").Value = 8.087130699999999
$ws.Range("P136").Value = 19.4532801
$ws.Range("Q136").Value = 14.4506336
$ws.Range("R136").Value = 13.53054
$ws.Range("S136").Value = 18.5071082
$ws.Range("T136").Value = 14.0283952
$ws.Range("U136").Value = 14.5827185
$ws.Range("V136").Value = 18.5492201
$ws.Range("W136").Value = 13.3782549
$ws.Range("X136").Value = 13.7620279
$ws.Range("Y136").Value = 9.234173200000001
$ws.Range("Z136").Value = 11.474599
$ws.Range("AA136").Value = 14.7731333
$ws.Range("AB136").Value = 13.9514174
$ws.Range("AD136").Value = 20.5963494
$ws.Range("AE136").Value = 9.8346175
$ws.Range("AF136").Value = 14.7292706
$ws.Range("AG136").Value = 16.7617261
$ws.Range("AH136").Value = 20.0518177
$ws.Range("AI136").Value = 11.0003244
$ws.Range("AJ136").Value = 13.9168783
$ws.Range("AK136").Value = 14.4731198
$ws.Range("AL136").Value = 13.0884216
$ws.Range("AM136").Value = 12.6239328
$ws.Range("AN136").Value = 13.6194837
$ws.Range("AO136").Value = 15.3055437
$ws.Range("AP136").Value = 12.02674
$ws.Range("AQ136").Value = 11.3726958
$ws.Range("AS136").Value = 11.7449412
$ws.Range("AT136").Value = 17.5838335
$ws.Range("AU136").Value = 20.8173489
$ws.Range("AV136").Value = 13.9945432
$ws.Range("AW136").Value = 16.4032226
$ws.Range("AX136").Value = 18.6003755
$ws.Range("AY136").Value = 13.4537752
$ws.Range("BA136").Value = 8.1673332
$ws.Range("BB136").Value = 13.2115063
$ws.Range("BC136").Value = 13.4848592
$ws.Range("BD136").Value = 13.073168
$ws.Range("BE136").Value = 13.0956232

# Row 137
$ws.Range("A137").Value = "15 06 2020"

# Row 138
$ws.Range("A138").Value = "16 06 2020"
